$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: row number -> values for columns B..N (13 columns)
# The underlying source values had been mis-parsed with "." used as a
# thousands separator (pt-BR style) instead of a decimal point, so the
# true integer magnitude is recovered here and totals (col N) are the
# recomputed sum of the corrected monthly values (B..M).
$rowData = @{}
$rowData[2] = @(317183,301322,317174,276684,276721,303261,303903,297064,300749,304548,288071,332637,3619317)
$rowData[3] = @(11783,10175,12628,1108,10858,1292,12996,9697,8713,9254,7883,7443,103830)
$rowData[4] = @(32876,34225,36639,29838,28588,29419,33996,29855,31316,31237,25779,2776,346544)
$rowData[5] = @(34419,34304,38879,34545,33554,35958,34066,27748,28413,28761,27706,29169,387522)
$rowData[6] = @(1605,1411,18095,17825,19127,17875,22443,18649,23443,18184,12851,14864,186372)
$rowData[7] = @(12983,19365,17045,17199,16207,18534,2049,14844,14083,15415,13485,13779,174988)
$rowData[8] = @(17766,8768,19418,19156,17993,20414,19061,18753,22086,229,32003,38609,234256)
$rowData[9] = @(4527,6423,4381,3663,387,485,6232,637,7762,7598,6998,8874,57967)
$rowData[10] = @(8425,7271,8039,8058,7459,7894,7986,7689,7811,793,7283,862,79570)
$rowData[11] = @(6562,7013,718,62,6578,7144,6561,5951,5645,5859,5277,5836,63206)
$rowData[12] = @(1249,1539,1933,1682,1896,2323,2957,2699,2653,3029,1465,138,23563)
$rowData[13] = @(5236,5139,5959,532,4532,5679,6123,5126,5078,5051,5446,6174,60075)
$rowData[14] = @(2711,2589,2727,2488,3017,3292,3476,3395,3251,3308,3281,3334,36869)
$rowData[15] = @(1629,139,1175,1662,138,1719,1886,1525,1688,1619,1248,1641,16069)
$rowData[16] = @(1321,1574,1852,1289,9510,1195,1409,1407,1419,1218,8530,9430,40154)
$rowData[17] = @(7290,6910,7340,5400,5970,7100,8230,8030,8240,8610,8060,8370,89550)
$rowData[18] = @(0,0,0,0,0,0,0,0,0,0,0,0,0)
$rowData[19] = @(0,0,0,0,0,0,0,0,0,0,0,0,0)
$rowData[20] = @(0,0,0,0,0,0,0,0,0,0,0,0,0)
$rowData[21] = @(100,80,150,0,120,0,0,0,0,0,100,0,550)
$rowData[22] = @(2868,3152,2939,1871,2302,3953,4745,5243,5137,5122,5404,6389,49125)
$rowData[23] = @(20,10,30,30,20,10,10,30,10,10,20,20,220)
$rowData[24] = @(4520,5180,3500,4130,3770,5110,4600,5540,4420,5540,5290,4930,56530)
$rowData[25] = @(0,0,0,0,0,0,0,0,0,0,0,0,0)
$rowData[26] = @(0,0,0,0,0,0,0,0,0,0,3530,3220,6750)
$rowData[27] = @(290,10,110,0,40,20,610,520,650,240,0,140,2630)
$rowData[28] = @(0,0,0,0,0,610,470,3940,400,180,1040,170,6810)
$rowData[29] = @(0,0,10,0,0,0,0,0,0,0,0,0,10)
$rowData[30] = @(5466,1056,1005,2829,2931,9486,11066,10261,864,9572,10718,9182,74436)
$rowData[31] = @(0,0,100,0,0,0,0,0,0,0,0,0,100)
$rowData[32] = @(0,0,40,0,0,70,0,0,0,0,0,0,110)
$rowData[33] = @(0,0,0,0,0,0,0,1244,30,1960,1370,1200,5804)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
